$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (row 9) down into the two new rows
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E11").PasteSpecial(-4122)

# Row 10: Change Number 08
$ws.Range("A10").Value = 41720
$ws.Range("B10").Value = "08"
$ws.Range("C10").Value = "JEB"
$ws.Range("D10").Value = "Initial Requirements matrix for Homework 3"
$ws.Range("E10").Value = "Done"

# Row 11: Change Number 09
$ws.Range("A11").Value = 41720
$ws.Range("B11").Value = "09"
$ws.Range("C11").Value = "JMR"
$ws.Range("D11").Value = "Project folder for homework 3"
$ws.Range("E11").Value = "Done"

# Update the active selection to match the new last cell
$ws.Range("A11").Select()
